$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BoM")

# Helper: write a literal text value into a cell while preserving the
# cell's existing style/number-format (plain .Value assignment on a
# numeric-looking string would get auto-parsed into a real number and
# also tends to drop trailing zeros / introduce float rounding, and
# explicitly touching .NumberFormat stamps a brand-new style index).
# Routing the text through a text-producing formula and then doing a
# Copy + PasteSpecial(values-only) converts it back to a literal string
# cell without disturbing the style index.
function Set-TextValue($range, $text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
    $excel.CutCopyMode = $false
}

# Date: regenerated BoM was produced a few hours later
$ws.Range("D5").Value = "2023-05-26_17-21-17"

# Footprint X for C4 (0.22uF cap)
Set-TextValue $ws.Range("O10") "2.9000"

# Footprint Y for C3 (10uF cap)
Set-TextValue $ws.Range("P12") "-11.1000"

# J1 connector: datasheet link replaced with placeholder, footprint X updated
$ws.Range("L14").Value = "~"
Set-TextValue $ws.Range("O14") "16.5000"

# L1 inductor: footprint X/Y updated
Set-TextValue $ws.Range("O15") "7.9000"
Set-TextValue $ws.Range("P15") "-7.1000"

# U2 (LT3494) footprint + footprint lib changed
Set-TextValue $ws.Range("G18") "DFN-8-1EP_2x3mm_P0.5mm_EP0.61x2.2mm"
Set-TextValue $ws.Range("H18") "Package_DFN_QFN"

# Swap the "changed field" highlight color between L14 (no longer flagged)
# and H18 (now flagged, since it picked up a real value)
$ws.Range("L14").Interior.Color = 9079551
$ws.Range("H18").Interior.Color = 12447999
